# The "Team ID" value cell currently reads "LTVIP2026TMIDS" + "64" + "787"
# (three runs). The edit changes the final run's text from "787" to "78",
# and adds a brand-new run containing "9" right after it - i.e. the last
# digit "7" is replaced by "9", but typed as a fresh run rather than simply
# rewriting the existing run's text.
$d = $word.ActiveDocument

# Locate the "787" run.
$run787 = $d.Content
$run787.Find.Execute("787", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Force Word to split "787" into two runs ("78" | "7") by nudging the
# formatting of the first two characters (color off, then back to the
# original value) without touching the text itself.
$firstTwo = $d.Range($run787.Start, $run787.Start + 2)
$firstTwo.Font.Color = 255
$firstTwo.Font.Color = 2236962

# Re-find "787" (Start/End are unchanged by the formatting nudge) and grab
# the now-isolated trailing "7" character as its own run.
$run787b = $d.Content
$run787b.Find.Execute("787", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lastChar = $d.Range($run787b.End - 1, $run787b.End)

# Nudge its formatting too, so replacing its text doesn't get folded back
# into the neighboring "78" run, then type the new digit.
$lastChar.Font.Color = 999999
$lastChar.Text = "9"

# Finally restore the correct font color on the freshly typed "9" run.
$lastChar.Font.Color = 2236962
